$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("evaluation")

# Update the date-range label in B14 (max drawdown duration period)
$ws.Range("B14").Value = "20171110-20190906"

# Update recalculated NAV-derived metrics in column B (rows 2-17)
$ws.Range("B2").Value  = -0.1764590984603687
$ws.Range("B3").Value  = 1.893044957579559
$ws.Range("B4").Value  = -3.503022884605334
$ws.Range("B5").Value  = 8.467771253307955
$ws.Range("B6").Value  = 13.13545815524904
$ws.Range("B7").Value  = 7.289169244858074
$ws.Range("B8").Value  = 28.09708728017304
$ws.Range("B9").Value  = 4.363555232530292
$ws.Range("B10").Value = 4.808098398429319
$ws.Range("B11").Value = 6.583639484494391
$ws.Range("B13").Value = 665
$ws.Range("B15").Value = 0.5955567959645179
$ws.Range("B16").Value = 0.4349435817249994
$ws.Range("B17").Value = 8.983719761783686
